$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "59.420.08"
$ws.Range("E2").Value = "  -5.67%  "

# Row 3
$ws.Range("D3").Value = "2.445.64"
$ws.Range("E3").Value = "  -8.86%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.13"
$ws.Range("E5").Value = "  -2.75%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.57"
$ws.Range("E6").Value = "  -6.32%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.25%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("E8").Value = "  -3.30%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0992"
$ws.Range("E9").Value = "  -6.01%  "

# Row 10
$ws.Range("E10").Value = "  -2.58%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.39"
$ws.Range("E11").Value = "  +5.53%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.350"
$ws.Range("E12").Value = "  -4.92%  "

# Row 13
$ws.Range("D13").Value = "2.878.30"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.16"
$ws.Range("E14").Value = "  -7.41%  "

# Row 15
$ws.Range("D15").Value = "59.372.38"
$ws.Range("E15").Value = "  -5.57%  "

# Row 16
$ws.Range("E16").Value = "  -5.87%  "

# Row 17
$ws.Range("D17").Value = "2.492.86"
$ws.Range("E17").Value = "  -7.20%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.17"
$ws.Range("E18").Value = "  -6.17%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.36"
$ws.Range("E19").Value = "  -4.65%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "324.88"
$ws.Range("E20").Value = "  -5.36%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.970"
$ws.Range("E21").Value = "  -2.87%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.71"
$ws.Range("E22").Value = "  -9.75%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.467"
$ws.Range("E23").Value = "  -7.34%  "

# Row 24
$ws.Range("E24").Value = "  -4.93%  "

# Row 25
$ws.Range("E25").Value = "  -3.79%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.976"
$ws.Range("E26").Value = "  -2.40%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.74"
$ws.Range("E27").Value = "  -4.79%  "

# Row 28
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.30"
$ws.Range("E28").Value = "  -2.70%  "

# Row 29
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.89"
$ws.Range("E29").Value = "  -1.94%  "

# Row 30
$ws.Range("E30").Value = "  -5.62%  "

# Row 31
$ws.Range("D31").Value = "0.0₃0770"
$ws.Range("E31").Value = "  -10.08%  "

# Row 32
$ws.Range("E32").Value = "  -0.11%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "156.20"
$ws.Range("E33").Value = "  -6.58%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.54"
$ws.Range("E34").Value = "  -5.59%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.34"
$ws.Range("E35").Value = "  -6.05%  "

# Row 36
$ws.Range("E36").Value = "  -5.57%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.74"
$ws.Range("E37").Value = "  -1.64%  "

# Row 38
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "312.86"
$ws.Range("E38").Value = "  -7.79%  "

# Row 39
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.76"
$ws.Range("E39").Value = "  -6.88%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.854"
$ws.Range("E40").Value = "  -8.17%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.76"
$ws.Range("E41").Value = "  -3.87%  "

# Row 42
$ws.Range("E42").Value = "  -5.77%  "

# Row 43
$ws.Range("E43").Value = "  -0.41%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.70"
$ws.Range("E44").Value = "  -3.10%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.582"
$ws.Range("E45").Value = "  -5.65%  "

# Row 46
$ws.Range("E46").Value = "  -3.41%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0524"
$ws.Range("E47").Value = "  -6.67%  "

# Row 48
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0230"
$ws.Range("E48").Value = "  -4.21%  "

# Row 49
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.04"
$ws.Range("E49").Value = "  -8.26%  "

# Row 50
$ws.Range("E50").Value = "  -8.74%  "

# Row 51
$ws.Range("D51").Value = "1.988.69"
$ws.Range("E51").Value = "  -4.87%  "
